$d = $word.ActiveDocument
$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# Find-ParagraphByText returns the Nth (1-based, default 1) paragraph whose
# text equals $oldText (ignoring the trailing paragraph-mark character).
function Find-ParagraphByText($oldText, $occurrence) {
    if (-not $occurrence) { $occurrence = 1 }
    $seen = 0
    foreach ($p in $d.Paragraphs) {
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $oldText) {
            $seen = $seen + 1
            if ($seen -eq $occurrence) {
                return $p
            }
        }
    }
    return $null
}

# Replace-ParagraphText rebuilds the whole paragraph from fresh XML so that
# structural bits (<w:pPr>, run-level <w:rPr>) survive the edit exactly as
# they were, instead of being silently dropped by a plain text/Find-Replace
# rewrite. InsertXML on a Range spanning a full paragraph (incl. its mark)
# automatically re-synthesizes any leading empty <w:r/> that the original
# paragraph had (tied to the paragraph-mark run props) - so we must NOT add
# one ourselves, or it gets duplicated.
function Replace-ParagraphText($oldText, $newText, $pPrXml, $rPrXml, $occurrence) {
    $p = Find-ParagraphByText $oldText $occurrence
    if ($p -eq $null) {
        throw "Paragraph not found for: $oldText"
    }
    $r = $d.Range($p.Range.Start, $p.Range.End)

    $inner = ""
    if ($pPrXml -ne "") {
        $inner += $pPrXml
    }
    if ($rPrXml -ne "") {
        $inner += "<w:r>$rPrXml<w:t>$newText</w:t></w:r>"
    } else {
        $inner += "<w:r><w:t>$newText</w:t></w:r>"
    }

    $xml = "<w:p $ns>$inner</w:p>"
    $r.InsertXML($xml)
}

$bulletPPr = "<w:pPr><w:pStyle w:val='ListBullet'/><w:spacing w:line='240' w:lineRule='auto'/><w:ind w:left='720'/></w:pPr>"

# Title (Heading1, top of document) - no leading empty run, no rPr override.
Replace-ParagraphText "Play Majestic Megaways Free - Review of iSoftBet's Slot Machine" `
    "Play Majestic Megaways and Win Big for Free" `
    "<w:pPr><w:pStyle w:val='Heading1'/></w:pPr>" "" 1

# "What we like" bullet points - ListBullet paragraphs each with a leading empty run.
Replace-ParagraphText "Megaways game engine provides up to 117,649 ways of winning" `
    "Megaways game engine offers up to 117,649 ways to win" $bulletPPr "" 1

Replace-ParagraphText "Free Spin bonus round with unlimited multiplier" `
    "Unique free spin bonus round with increasing multipliers" $bulletPPr "" 1

Replace-ParagraphText "Full Moon Mystery symbol adds more excitement to the game" `
    "Full Moon Mystery symbol provides additional win potential" $bulletPPr "" 1

Replace-ParagraphText "Impressive graphics and immersive audio effects enhance the gaming experience" `
    "Engaging gameplay and impressive graphics" $bulletPPr "" 1

# "What we don't like" bullet points
Replace-ParagraphText "Does not offer progressive jackpot" `
    "Limited variety of bonus features" $bulletPPr "" 1

Replace-ParagraphText "No gamble feature available" `
    "May be complex for new slot players" $bulletPPr "" 1

# Bold title repeated near the bottom of the document - leading empty run, bold rPr.
# This is the 1st (and now only) remaining paragraph with this exact text,
# since the Heading1 occurrence above was already replaced.
Replace-ParagraphText "Play Majestic Megaways Free - Review of iSoftBet's Slot Machine" `
    "Play Majestic Megaways and Win Big for Free" "" "<w:rPr><w:b/></w:rPr>" 1

# Italic meta description - leading empty run, italic rPr.
Replace-ParagraphText "Read our review of iSoftBet's Majestic Megaways slot machine. Play for free and enjoy the exclusive Megaways game engine, Free Spins bonus round, and Full Moon Mystery symbol." `
    "Read our review of Majestic Megaways and discover its exciting features. Play for free and win big!" "" "<w:rPr><w:i/></w:rPr>" 1
